# "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet "dados" held a table split into two labelled sub-blocks, each
# introduced by a label-only row in column A:
#   row 5 -> "situação do domicílio"               (blank B:F)
#   row 8 -> "grandes regiões e unidades da federação" (blank B:F)
# Those section-header rows are removed outright (not just their text),
# so every row below shifts up and the "urbana"/"rural" and
# region/UF rows become contiguous with the "brasil" row above them.
# Also corrects the row-2 column headers, which were auto-generated
# "unnamed: N_level_1" placeholders that should read "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Delete bottom row first so the row number of the one still pending
# deletion (row 5) doesn't shift.
$ws.Range("A8").EntireRow.Delete()
$ws.Range("A5").EntireRow.Delete()

# Fix the row-2 placeholder headers.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
